# Update cryptos list on worksheet with latest prices / volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    # Force the value to be stored as text (avoids Excel auto-converting
    # numeric-looking strings like "230.61" into real numbers), while
    # keeping the cell style identical to an unstyled / default cell.
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "93.690.24"
$ws.Range("E2").Value = "  -1.52%  "

# Row 3 - Ethereum
Set-TextCell "D3" "3.322.38"
$ws.Range("E3").Value = "  -3.53%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - Solana
Set-TextCell "D5" "230.61"
$ws.Range("E5").Value = "  -3.45%  "

# Row 6 - BNB
Set-TextCell "D6" "617.03"
$ws.Range("E6").Value = "  -3.82%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -4.42%  "

# Row 8 - Dogecoin
$ws.Range("E8").Value = "  -3.49%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  -0.02%  "

# Row 10 - Cardano
Set-TextCell "D10" "0.924"
$ws.Range("E10").Value = "  -6.42%  "

# Row 11 - LidoStakedEther
Set-TextCell "D11" "3.319.24"
$ws.Range("E11").Value = "  -3.68%  "

# Row 12 - Avalanche
Set-TextCell "D12" "41.98"
$ws.Range("E12").Value = "  +1.09%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -1.92%  "

# Row 14 - WrappedBTC
Set-TextCell "D14" "93.573.41"
$ws.Range("E14").Value = "  -1.35%  "

# Row 15 - Toncoin
Set-TextCell "D15" "5.93"
$ws.Range("E15").Value = "  -2.42%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextCell "D16" "3.942.41"
$ws.Range("E16").Value = "  -3.42%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  -4.91%  "

# Row 18 - Polkadot
Set-TextCell "D18" "8.07"
$ws.Range("E18").Value = "  -3.93%  "

# Row 19 - WrappedEther
Set-TextCell "D19" "3.322.73"
$ws.Range("E19").Value = "  -4.08%  "

# Row 20 - Chainlink
Set-TextCell "D20" "17.13"
$ws.Range("E20").Value = "  -4.57%  "

# Row 21 - Uniswap
Set-TextCell "D21" "10.87"
$ws.Range("E21").Value = "  -6.04%  "

# Row 22 - SuiNetwork
$ws.Range("E22").Value = "  +9.61%  "

# Row 23 - BitcoinCash
Set-TextCell "D23" "493.87"
$ws.Range("E23").Value = "  -1.49%  "

# Row 24 - Stellar
Set-TextCell "D24" "0.445"
$ws.Range("E24").Value = "  -12.84%  "

# Row 25
$ws.Range("E25").Value = "  -5.06%  "

# Row 26
Set-TextCell "D26" "6.03"
$ws.Range("E26").Value = "  -6.27%  "

# Row 27
Set-TextCell "D27" "89.63"
$ws.Range("E27").Value = "  -1.89%  "

# Row 28 - Aptos
Set-TextCell "D28" "11.72"
$ws.Range("E28").Value = "  -2.40%  "

# Row 29 - WrappedeETH
Set-TextCell "D29" "3.504.23"

# Row 30 - Dai
$ws.Range("E30").Value = "  +0.18%  "

# Row 31 & 32 - swap InternetComputer(DFINITY) and Hedera
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D31" "0.139"
$ws.Range("E31").Value = "  +1.66%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D32" "11.01"
$ws.Range("E32").Value = "  -5.70%  "

# Row 33 - PancakeSwap
Set-TextCell "D33" "2.61"
$ws.Range("E33").Value = "  -4.18%  "

# Row 34 - Binance-PegBSC-USD
Set-TextCell "D34" "0.990"
$ws.Range("E34").Value = "  -1.10%  "

# Row 35 - Cronos
Set-TextCell "D35" "0.173"
$ws.Range("E35").Value = "  -4.95%  "

# Row 36 - EthereumClassic
Set-TextCell "D36" "28.23"
$ws.Range("E36").Value = "  -8.22%  "

# Row 37 - PolygonEcosystemToken
$ws.Range("E37").Value = "  -6.67%  "

# Row 38 - Bittensor
Set-TextCell "D38" "527.96"
$ws.Range("E38").Value = "  +3.15%  "

# Row 39 & 40 - swap RenderToken and USDe
$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell "D39" "1.00"
$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextCell "D40" "7.37"
$ws.Range("E40").Value = "  -4.39%  "

# Row 41 - Kaspa
Set-TextCell "D41" "0.147"
$ws.Range("E41").Value = "  -1.53%  "

# Row 42 - Fetch.AI
$ws.Range("E42").Value = "  -5.36%  "

# Row 43 - ARBITRUM
Set-TextCell "D43" "0.858"
$ws.Range("E43").Value = "  -5.69%  "

# Row 44 - WhiteBITCoin
Set-TextCell "D44" "24.05"
$ws.Range("E44").Value = "  -0.30%  "

# Row 45 - MantraDAO
Set-TextCell "D45" "3.65"
$ws.Range("E45").Value = "  +4.13%  "

# Row 46 - VeChain
$ws.Range("E46").Value = "  +0.11%  "

# Row 47 - ImmutableX
$ws.Range("E47").Value = "  -1.27%  "

# Row 48 - Filecoin
Set-TextCell "D48" "5.37"
$ws.Range("E48").Value = "  -2.84%  "

# Row 49 & 50 - swap OKB and Stacks
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell "D49" "2.10"
$ws.Range("E49").Value = "  -1.55%  "

$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D50" "52.58"
$ws.Range("E50").Value = "  -1.72%  "

# Row 51 - Cosmos
Set-TextCell "D51" "7.92"
$ws.Range("E51").Value = "  -0.34%  "
